# Update "epexspot_prices.xlsx" with the newest day of data:
#  - Prix Spot: new column BU ("25-aug") with 24 hourly prices
#  - Gaz:       two new rows (2025-08-23 / 2025-08-24), carrying the last price forward
#  - CO2:       two new rows (2025-08-23 / 2025-08-24), price not yet published (blank)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot" -> add column BU (25-aug)
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, formatted like the preceding day header (BT1)
$wsPrix.Cells.Item(1, 73).Value = "25-aug"
$wsPrix.Range("BT1").Copy()
$wsPrix.Cells.Item(1, 73).PasteSpecial(-4122)

$prixValues = @(
    95.28,
    87.03,
    74.55,
    60.88,
    62.43,
    80.16,
    91.90000000000001,
    102.52,
    101.82,
    91.97,
    73.5,
    35.02,
    10,
    5.93,
    6.78,
    25.2,
    56.43,
    76.19,
    97.43000000000001,
    114.78,
    125.95,
    126.28,
    115.04,
    105.15
)

$row = 2
foreach ($val in $prixValues) {
    $wsPrix.Cells.Item($row, 73).Value = $val
    $row++
}

# ---------------------------------------------------------------------------
# Sheet "Gaz" -> append the two missing dates (last known price carried fwd)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$gazDates = @("2025-08-23", "2025-08-24")
$gazValues = @(32.2, 32.2)

$row = 70
for ($i = 0; $i -lt $gazDates.Length; $i++) {
    # leading apostrophe forces text (no auto date conversion), then copy the
    # neighbouring cell's format so no stray "quote prefix" style is left behind
    $wsGaz.Cells.Item($row, 1).Value = "'" + $gazDates[$i]
    $wsGaz.Range("A69").Copy()
    $wsGaz.Cells.Item($row, 1).PasteSpecial(-4122)

    $wsGaz.Cells.Item($row, 2).Value = $gazValues[$i]
    $row++
}

# ---------------------------------------------------------------------------
# Sheet "CO2" -> append the two missing dates (price not published yet)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$co2Dates = @("2025-08-23", "2025-08-24")

$row = 70
for ($i = 0; $i -lt $co2Dates.Length; $i++) {
    $wsCo2.Cells.Item($row, 1).Value = "'" + $co2Dates[$i]
    $wsCo2.Range("A69").Copy()
    $wsCo2.Cells.Item($row, 1).PasteSpecial(-4122)

    # B column stays blank for this date, same as row 69 (empty text cell)
    $wsCo2.Cells.Item($row, 2).Value = "'"
    $wsCo2.Range("B69").Copy()
    $wsCo2.Cells.Item($row, 2).PasteSpecial(-4122)

    $row++
}

Write-Host "done"
